# ==========================================================================
# Update datamodel and mvp config example
# Applies the structural + value changes described in the commit diff to
# the five worksheets of the workbook: gridconnections, gridnodes (no
# change), actors, policies, contracts.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# Sheet: gridconnections
# --------------------------------------------------------------------
$ws = $wb.Worksheets.Item("gridconnections")

# Insert new header columns J:N (charging_mode, battery_mode,
# nfATO_capacity_kw, nfATO_starttime, nfATO_endtime) ahead of the old
# J/K headers (insulation_label / heating_type), which shift right to O/P.
# Copy the bold/bordered header style from the existing I1 header cell
# across the whole newly used header range first, so every new header
# cell matches the look of the existing ones.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:P1").PasteSpecial(-4122) | Out-Null

$ws.Range("J1").Value = "charging_mode"
$ws.Range("K1").Value = "battery_mode"
$ws.Range("L1").Value = "nfATO_capacity_kw"
$ws.Range("M1").Value = "nfATO_starttime"
$ws.Range("N1").Value = "nfATO_endtime"
$ws.Range("O1").Value = "insulation_label"
$ws.Range("P1").Value = "heating_type"

# Row 2 (BUILDING / LOGISTICS / b1)
$ws.Range("F2").Value = 750
$ws.Range("I2").Value = "['EHGV', 'EHGV', 'EHGV', 'EHGV', 'EHGV', 'EHGV', 'Diesel_Truck', 'Grid_battery_7MWh', 'Building_solarpanels_0kWp', 'Building_gas_burner_60kW', 'Building_solarpanels_10kWp']"
$ws.Range("J2").Value = "MAX_POWER"
$ws.Range("K2").Value = "BALANCE"
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = "NONE"
$ws.Range("P2").Value = "GASBURNER"

# Row 3 (INDUSTRY / INDUSTRY_OTHER / b2)
$ws.Range("I3").Value = "['INDUSTRY_OTHER_HEAT_DEMAND', 'Building_solarpanels_0kWp', 'Building_gas_burner_60kW']"
# old heating_type value ("GASBURNER") lived in K3; that column is now a
# blank inserted column, so clear it before writing the shifted value to P3
$ws.Range("K3").Value = ""
$ws.Range("P3").Value = "GASBURNER"

# Row 5 (GRIDBATTERY / b4) gains a battery_mode value
$ws.Range("K5").Value = "BALANCE"

# --------------------------------------------------------------------
# Sheet: actors
# --------------------------------------------------------------------
$ws = $wb.Worksheets.Item("actors")

# New header columns F:H (nfATO_capacity_kw, nfATO_starttime, nfATO_endtime)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

$ws.Range("F1").Value = "nfATO_capacity_kw"
$ws.Range("G1").Value = "nfATO_starttime"
$ws.Range("H1").Value = "nfATO_endtime"

# Row 7 (sup1 actor) becomes a NonFirmActor with nfATO values
$ws.Range("C7").Value = "NonFirmActor"
$ws.Range("F7").Value = 900
$ws.Range("G7").Value = 18
$ws.Range("H7").Value = 7

# --------------------------------------------------------------------
# Sheet: policies
# --------------------------------------------------------------------
$ws = $wb.Worksheets.Item("policies")

# Existing "value"/"unit" columns are stored as text, even when the
# content looks numeric - force text on every cell in column E (and the
# newly written column F values) using a leading apostrophe so Excel
# does not silently convert them to numbers/booleans.
$ws.Range("E3").Value = "'0"
$ws.Range("E4").Value = "'0.3"
$ws.Range("F4").Value = "eurpkWh"
$ws.Range("E5").Value = "'0.5"

# New policy rows 8-15
$ws.Range("C8").Value = "Policy"
$ws.Range("D8").Value = "Fixed_electricity_price"
$ws.Range("E8").Value = "'0.21"
$ws.Range("F8").Value = "EUR p kWh"
$ws.Range("G8").Value = "Fixed_electricity_price"

$ws.Range("C9").Value = "Policy"
$ws.Range("D9").Value = "Fixed_heat_price"
$ws.Range("E9").Value = "'0.10"
$ws.Range("F9").Value = "EUR p kWh"
$ws.Range("G9").Value = "Fixed_heat_price"

$ws.Range("C10").Value = "Policy"
$ws.Range("D10").Value = "Fixed_methane_price"
$ws.Range("E10").Value = "'0.05"
$ws.Range("F10").Value = "EUR p kWh"
$ws.Range("G10").Value = "Fixed_methane_price"

$ws.Range("C11").Value = "Policy"
$ws.Range("D11").Value = "Fixed_hydrogen_price"
$ws.Range("E11").Value = "'0.3"
$ws.Range("F11").Value = "EUR p kWh"
$ws.Range("G11").Value = "Fixed_hydrogen_price"

$ws.Range("C12").Value = "Policy"
$ws.Range("D12").Value = "Energy_supplier_electricity_price_margin"
$ws.Range("E12").Value = "'0.17"
$ws.Range("F12").Value = "fr"
$ws.Range("G12").Value = "p_variableElectricityPriceOverNational_eurpkWh"

$ws.Range("C13").Value = "Policy"
$ws.Range("D13").Value = "Fixed_electricity_feed_in_tariff"
$ws.Range("E13").Value = "'0.25"
$ws.Range("F13").Value = "EUR_pKWh"
$ws.Range("G13").Value = "p_fixedFeedinTariff_eurpkWh"

$ws.Range("C14").Value = "Policy"
$ws.Range("D14").Value = "Fixed_diesel_price"
$ws.Range("E14").Value = "'0.15"
$ws.Range("F14").Value = "EUR_pKWh"
$ws.Range("G14").Value = "Fixed_diesel_price"

$ws.Range("C15").Value = "Policy"
$ws.Range("D15").Value = "Time_buffer_for_spread_charging"
$ws.Range("E15").Value = "'60"
$ws.Range("F15").Value = "minutes"
$ws.Range("G15").Value = "Time_buffer_for_spread_charging, Integer value"

# --------------------------------------------------------------------
# Sheet: contracts
# --------------------------------------------------------------------
$ws = $wb.Worksheets.Item("contracts")

# Row 2 unchanged: CONNECTIONOWNER | DEFAULT | Contract | ENERGYHOLON | com1

# Row 3: DEFAULT/ENERGYHOLON/com2 -> VARIABLE/ENERGYSUPPLIER/com1
$ws.Range("B3").Value = "VARIABLE"
$ws.Range("D3").Value = "ENERGYSUPPLIER"
$ws.Range("E3").Value = "com1"

# Row 4: only actor_id changes com3 -> com2
$ws.Range("E4").Value = "com2"

# Row 5: DEFAULT/ENERGYHOLON/com4 -> VARIABLE/ENERGYSUPPLIER/com2
$ws.Range("B5").Value = "VARIABLE"
$ws.Range("D5").Value = "ENERGYSUPPLIER"
$ws.Range("E5").Value = "com2"

# Row 6 (was ENERGYHOLON/GOPACS/Contract/GRIDOPERATOR/hol1) is overwritten
# with a new CONNECTIONOWNER/DEFAULT/Contract/ENERGYHOLON/com3 row
$ws.Range("A6").Value = "CONNECTIONOWNER"
$ws.Range("B6").Value = "DEFAULT"
$ws.Range("C6").Value = "Contract"
$ws.Range("D6").Value = "ENERGYHOLON"
$ws.Range("E6").Value = "com3"

# Row 7 (new): CONNECTIONOWNER/DEFAULT/Contract/ENERGYHOLON/com4
$ws.Range("A7").Value = "CONNECTIONOWNER"
$ws.Range("B7").Value = "DEFAULT"
$ws.Range("C7").Value = "Contract"
$ws.Range("D7").Value = "ENERGYHOLON"
$ws.Range("E7").Value = "com4"

# Row 8 (new): the original GOPACS/GRIDOPERATOR row, moved down from row 6
$ws.Range("A8").Value = "ENERGYHOLON"
$ws.Range("B8").Value = "GOPACS"
$ws.Range("C8").Value = "Contract"
$ws.Range("D8").Value = "GRIDOPERATOR"
$ws.Range("E8").Value = "hol1"

Write-Host "edit.ps1 completed"
